# Update temperature (column D) and uncertainty (column E) values for rows 2-19
# on the active worksheet, reflecting refreshed bootstrap results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = 15.0677871054159;  E = 0.8086528589963745 },
    @{ Row = 3;  D = 13.65273615994774; E = 2.056412934489027 },
    @{ Row = 4;  D = 15.84884973754504; E = 0.4737798588840921 },
    @{ Row = 5;  D = 14.82296216149963; E = 1.470271652439198 },
    @{ Row = 6;  D = 16.27049478550594; E = 0.5175175310755836 },
    @{ Row = 7;  D = 15.36140106242473; E = 1.160568058301429 },
    @{ Row = 8;  D = 16.27619821733979; E = 0.4962008151841023 },
    @{ Row = 9;  D = 15.48154308277698; E = 0.7801327716119801 },
    @{ Row = 10; D = 16.82246103450655; E = 0.521185082417052 },
    @{ Row = 11; D = 16.13016089284572; E = 0.8896458082423097 },
    @{ Row = 12; D = 20.04140488260402; E = 2.00159589401732 },
    @{ Row = 13; D = 17.34250682055415; E = 0.5344189831229733 },
    @{ Row = 14; D = 20.93593117027743; E = 2.659476564274648 },
    @{ Row = 15; D = 17.63861984181357; E = 1.00072081656106 },
    @{ Row = 16; D = 21.80860802039622; E = 2.621958038061083 },
    @{ Row = 17; D = 17.41312086157644; E = 1.098932289689143 },
    @{ Row = 18; D = 22.07572639411203; E = 2.706390963726625 },
    @{ Row = 19; D = 18.18954593210566; E = 1.014170911387849 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
